$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New wallet-ledger transactions recorded through 14-Feb (rows 53-60).
# Column layout: A=Date, B=Amt Debited, C=Amt Credited, D=Transaction Mode, E=Wallet Balance (running formula)

$dateFmt = "[$-409]d\-mmm\-yyyy;@"
$balanceFormula = '=IF(A{0}="","",SUM(E{1}-B{0}+C{0}))'

function Set-LedgerRow($row, $date, $debit, $credit, $mode) {
    $ws.Range("A$row").Value = $date
    $ws.Range("A$row").NumberFormat = $dateFmt

    if ($debit -ne $null) {
        $ws.Range("B$row").Value = $debit
    }
    if ($credit -ne $null) {
        $ws.Range("C$row").Value = $credit
    }

    $ws.Range("D$row").Value = $mode

    $prevRow = $row - 1
    $ws.Range("E$row").Formula = [string]::Format($balanceFormula, $row, $prevRow)
}

Set-LedgerRow 53 43870 46280 $null "Ordered Amount"
Set-LedgerRow 54 43871 28080 $null "Ordered Amount"
Set-LedgerRow 55 43871 $null 111925 "Manual Added"
Set-LedgerRow 56 43872 42640 $null "Ordered Amount"
Set-LedgerRow 57 43873 58760 $null "Ordered Amount"
Set-LedgerRow 58 43874 42640 $null "Ordered Amount"
Set-LedgerRow 59 43875 47840 $null "Ordered Amount"
Set-LedgerRow 60 43875 $null 32733 "Manual Added"

# Update the last active-cell selection left behind by the editor (E57 -> E61)
[void]$ws.Range("E61").Select()
